# Update "想去人数" (column F) figures across the four worksheets to the
# values generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 2406
        7  = 264
        11 = 1042
        12 = 839
        13 = 89
        14 = 835
        15 = 1451
        16 = 689
        18 = 36
        19 = 340
    }
    "演出" = @{
        19 = 149
        24 = 86
        26 = 43
        38 = 330
        46 = 296
    }
    "本地生活" = @{
        4  = 2478
        6  = 2488
        12 = 349
        13 = 2757
        15 = 655
    }
    "全部类型" = @{
        2  = 2478
        7  = 2757
        10 = 655
        17 = 264
        20 = 1042
        21 = 839
        22 = 89
        23 = 835
        28 = 689
        32 = 340
        33 = 86
        34 = 43
        43 = 330
        48 = 296
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
